# fall 23 week 1 inputs
# Update a handful of matchup-average values on the "Nine" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D3").Value = 10.69
$ws.Range("G3").Value = 11.56

$ws.Range("C4").Value = 9.31
$ws.Range("E4").Value = 10.01
$ws.Range("F4").Value = 9.95

$ws.Range("D5").Value = 9.99
$ws.Range("F5").Value = 10.28
$ws.Range("H5").Value = 8.1

$ws.Range("D6").Value = 10.05
$ws.Range("E6").Value = 9.72
$ws.Range("G6").Value = 10.43

$ws.Range("C7").Value = 8.44
$ws.Range("F7").Value = 9.57

$ws.Range("E8").Value = 11.9
